$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-01 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-12-02 Saturday", 2) | Out-Null
$d.Content.Find.Execute("42÷4=10, 2", $true, $true, $false, $false, $false, $true, 1, $false, "79÷8=9, 7", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $true, $false, $false, $false, $true, 1, $false, "73÷3=24, 1", 2) | Out-Null
$d.Content.Find.Execute("14÷3=4, 2", $true, $true, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("45÷7=6, 3", $true, $true, $false, $false, $false, $true, 1, $false, "33÷2=16, 1", 2) | Out-Null
$d.Content.Find.Execute("88÷5=17, 3", $true, $true, $false, $false, $false, $true, 1, $false, "26÷6=4, 2", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $true, $false, $false, $false, $true, 1, $false, "45÷6=7, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷2=24, 0", $true, $true, $false, $false, $false, $true, 1, $false, "37÷8=4, 5", 2) | Out-Null
$d.Content.Find.Execute("38÷5=7, 3", $true, $true, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2) | Out-Null
$d.Content.Find.Execute("62÷7=8, 6", $true, $true, $false, $false, $false, $true, 1, $false, "41÷8=5, 1", 2) | Out-Null
$d.Content.Find.Execute("81÷8=10, 1", $true, $true, $false, $false, $false, $true, 1, $false, "36÷5=7, 1", 2) | Out-Null
$d.Content.Find.Execute("10÷9=1, 1", $true, $true, $false, $false, $false, $true, 1, $false, "99÷8=12, 3", 2) | Out-Null
$d.Content.Find.Execute("82÷4=20, 2", $true, $true, $false, $false, $false, $true, 1, $false, "61÷2=30, 1", 2) | Out-Null
$d.Content.Find.Execute("77÷9=8, 5", $true, $true, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷7=11, 1", $true, $true, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2) | Out-Null
$d.Content.Find.Execute("19÷9=2, 1", $true, $true, $false, $false, $false, $true, 1, $false, "88÷2=44, 0", 2) | Out-Null
$d.Content.Find.Execute("22÷2=11, 0", $true, $true, $false, $false, $false, $true, 1, $false, "65÷8=8, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷6=6, 3", $true, $true, $false, $false, $false, $true, 1, $false, "68÷8=8, 4", 2) | Out-Null
$d.Content.Find.Execute("35÷4=8, 3", $true, $true, $false, $false, $false, $true, 1, $false, "73÷7=10, 3", 2) | Out-Null
$d.Content.Find.Execute("66÷5=13, 1", $true, $true, $false, $false, $false, $true, 1, $false, "96÷2=48, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷3=18, 1", $true, $true, $false, $false, $false, $true, 1, $false, "80÷7=11, 3", 2) | Out-Null
$d.Content.Find.Execute("17÷8=2, 1", $true, $true, $false, $false, $false, $true, 1, $false, "99÷9=11, 0", 2) | Out-Null
$d.Content.Find.Execute("29÷6=4, 5", $true, $true, $false, $false, $false, $true, 1, $false, "51÷9=5, 6", 2) | Out-Null
$d.Content.Find.Execute("72÷8=9, 0", $true, $true, $false, $false, $false, $true, 1, $false, "20÷6=3, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷8=2, 5", $true, $true, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("50÷9=5, 5", $true, $true, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
